$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "Favored Team Final Score"
$ws.Range("I1").Value = "Underdog Team Final Score"
$ws.Range("L1").Value = "Home Team Cover"
$ws.Range("J1").Value = "Favored Team Cover?"
$ws.Range("K1").Value = "Underdog Team Cover?"
